$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.01
$ws.Range("G4").Value = 0.01
$ws.Range("F5").Value = 0.02
$ws.Range("G5").Value = 0.98
$ws.Range("H5").Value = 0
$ws.Range("F6").Value = 0.01
